$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.379.84"
$ws.Range("E2").Value = "  -2.52%  "
$ws.Range("D3").Value = "3.307.31"
$ws.Range("E3").Value = "  -3.18%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'558.17"
$ws.Range("E5").Value = "  -3.06%  "
$ws.Range("D6").Value = "'142.48"
$ws.Range("E6").Value = "  -4.13%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.307.60"
$ws.Range("E8").Value = "  -3.18%  "
$ws.Range("E9").Value = "  -2.99%  "
$ws.Range("D10").Value = "'7.86"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("E11").Value = "  -3.39%  "
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("D13").Value = "3.874.67"
$ws.Range("E13").Value = "  -3.10%  "
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("E15").Value = "  -5.12%  "
$ws.Range("D16").Value = "3.301.70"
$ws.Range("E16").Value = "  -3.02%  "
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").Value = "60.388.17"
$ws.Range("E18").Value = "  -2.50%  "
$ws.Range("D19").Value = "'6.11"
$ws.Range("E19").Value = "  -3.88%  "
$ws.Range("D20").Value = "'14.13"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("D22").Value = "'375.22"
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("D23").Value = "'74.54"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("E25").Value = "  -4.78%  "
$ws.Range("D26").Value = "3.446.17"
$ws.Range("E26").Value = "  -3.08%  "
$ws.Range("E27").Value = "  -7.63%  "
$ws.Range("E28").Value = "  -4.17%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -5.86%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  -3.45%  "
$ws.Range("E33").Value = "  -3.93%  "
$ws.Range("D34").Value = "'22.65"
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("E35").Value = "  -6.57%  "
$ws.Range("E36").Value = "  -6.03%  "
$ws.Range("D37").Value = "'166.64"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("E38").Value = "  -3.48%  "
$ws.Range("D39").Value = "'6.71"
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("D40").Value = "3.339.46"
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("D41").Value = "'26.83"
$ws.Range("E41").Value = "  -13.05%  "
$ws.Range("D42").Value = "'0.0728"
$ws.Range("E42").Value = "  -5.95%  "
$ws.Range("D43").Value = "'41.92"
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("E44").Value = "  -2.94%  "
$ws.Range("E45").Value = "  -4.83%  "
$ws.Range("E46").Value = "  -5.14%  "
$ws.Range("E47").Value = "  -3.69%  "
$ws.Range("D48").Value = "2.358.82"
$ws.Range("E48").Value = "  -7.23%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  -6.61%  "
$ws.Range("D51").Value = "'21.33"
$ws.Range("E51").Value = "  -4.63%  "
